$d = $word.ActiveDocument
$t = $d.Tables(1)
$cell = $t.Cell(1,1)
$r = $cell.Range
$txt = $r.Text
Write-Output ("before: [" + $txt + "]")
$r.Text = $txt
Write-Output "done"
